# Auto-update draw results: append the 2025-12-20 Pick 3 draw as a new row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 95
$rng = $ws.Range("A" + $newRow + ":E" + $newRow)

# The existing data rows store every value as text (dates, 6-digit phase
# codes, and the "D-D-D" result are all plain strings), even though some of
# them look like numbers/dates. Force the new cells to text formatting
# before writing so Excel doesn't auto-convert them into a date serial
# number / numeric value.
$rng.NumberFormat = "@"

$ws.Range("A" + $newRow).Value = "2025-12-20"
$ws.Range("B" + $newRow).Value = "Pick 3"
$ws.Range("C" + $newRow).Value = "251220"
$ws.Range("D" + $newRow).Value = "4-7-7"
$ws.Range("E" + $newRow).Value = "2025-12-20T21:37:21.337+04:00"

# Restore the default "Normal" style so the new cells don't end up carrying
# an explicit text-number-format style (matching the rest of the sheet,
# which relies on default/general formatting).
$rng.Style = "Normal"
